$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the login test-data rows (A2:B6) -----------------------------
# Column A holds the "username" (mailto-hyperlinked) values, column B the
# "password" values. We only touch the cell .Value so the pre-existing
# hyperlinks (which still point at the original mailto: targets) stay
# attached to the cells, exactly like the source edit.
$ws.Range("A2").Value = "arjavopencart@gmail.com"

$ws.Range("A3").Value = "test123@gmail.com"
$ws.Range("B3").Value = "test123"

$ws.Range("A4").Value = "test456@yahoo.com"
$ws.Range("B4").Value = "test456"

$ws.Range("A5").Value = "test@amazon.com"
$ws.Range("B5").Value = "test789"

$ws.Range("A6").Value = "arjav@gmail.com"
$ws.Range("B6").Value = "arjav"

# --- Formatting touch-ups that came along with the edit -------------------
# The username column shrinks from the oversized 16pt hyperlink look down to
# the normal 11pt hyperlink style.
$ws.Range("A2").Font.Size = 11
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)

# The old green "Valid" highlight on C6 is replaced with the plain style
# already used by the other result cells.
$ws.Range("C2").Copy()
$ws.Range("C6").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Selection / active cell ----------------------------------------------
$ws.Range("J5").Select()
